$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Task #4: Implement a function to delete a tag from all items containing
# the tag" picked up one more day of work (Day 2 -> 3 tasks done), per the
# commit message "got ALL delete button semi working".
$ws.Range("D7").Value = 3

# Leave the cursor where the edit happened.
$ws.Range("D7").Select()
